$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 54, pushing the existing
# rows 54-64 down to become rows 56-66 (with all their values/formatting
# carried along automatically by Excel).
$ws.Rows("54:55").Insert()

# New row 54: Especial, Volumen 20 (same as old row 54 but with updated
# Fecha and Volumen).
$ws.Cells.Item(54, 1).Value = 10
$ws.Cells.Item(54, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(54, 3).Value = "La Araucanía"
$ws.Cells.Item(54, 4).Value = 44474
$ws.Cells.Item(54, 5).Value = 9
$ws.Cells.Item(54, 6).Value = "Fruta"
$ws.Cells.Item(54, 7).Value = 100107
$ws.Cells.Item(54, 8).Value = "Otros"
$ws.Cells.Item(54, 9).Value = 100107002
$ws.Cells.Item(54, 10).Value = "Chirimoya"
$ws.Cells.Item(54, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(54, 12).Value = "Especial"
$ws.Cells.Item(54, 13).Value = 20
$ws.Cells.Item(54, 14).Value = 3500
$ws.Cells.Item(54, 15).Value = 3500
$ws.Cells.Item(54, 16).Value = 3500
$ws.Cells.Item(54, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(54, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(54, 19).Value = 3500
$ws.Cells.Item(54, 20).Value = 1

# New row 55: Primera, Volumen 60 (same as old row 55 but with updated
# Fecha and Volumen).
$ws.Cells.Item(55, 1).Value = 10
$ws.Cells.Item(55, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(55, 3).Value = "La Araucanía"
$ws.Cells.Item(55, 4).Value = 44474
$ws.Cells.Item(55, 5).Value = 9
$ws.Cells.Item(55, 6).Value = "Fruta"
$ws.Cells.Item(55, 7).Value = 100107
$ws.Cells.Item(55, 8).Value = "Otros"
$ws.Cells.Item(55, 9).Value = 100107002
$ws.Cells.Item(55, 10).Value = "Chirimoya"
$ws.Cells.Item(55, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(55, 12).Value = "Primera"
$ws.Cells.Item(55, 13).Value = 60
$ws.Cells.Item(55, 14).Value = 3000
$ws.Cells.Item(55, 15).Value = 3000
$ws.Cells.Item(55, 16).Value = 3000
$ws.Cells.Item(55, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(55, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(55, 19).Value = 3000
$ws.Cells.Item(55, 20).Value = 1
